# Edit script: update Saldo.xlsx "Export" sheet per commit diff.
#
# Summary of changes:
#  1) Rows for account 005002390 (LUCIANO, 11000) and 004332783 (IRON, 6000)
#     — originally rows 4 and 5 — are removed and replaced by four new rows:
#       004479463 HENRIQUE 47959.99
#       004240014 ISABELE  25461.1
#       004752461 SERGIO   6623.66
#       004752494 SERGIO   6623.66
#  2) Four individual rows are removed elsewhere in the sheet:
#       004767746 ISABELE  77.36
#       004484207 FLAVIA   64.04
#       004497875 HENRIQUE 60.36
#       004752494 SERGIO   2.52

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: remove the four single rows, working bottom-to-top so that
#     earlier row numbers stay valid while later ones are deleted. These
#     all sit below row 5, so they do not disturb the rows-4/5 replacement
#     done afterwards.

$singleRowTargets = @("004767746", "004484207", "004497875", "004752494")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

foreach ($acct in $singleRowTargets) {
    for ($r = $lastRow; $r -ge 2; $r--) {
        $cellVal = $ws.Cells.Item($r, 1).Value()
        if ($cellVal -eq $acct) {
            $ws.Rows.Item($r).Delete()
            break
        }
    }
}

# --- Step 2: replace the LUCIANO / IRON rows (rows 4 and 5) with the four
#     new rows.

# Delete the two old rows (IRON first so LUCIANO's row index 4 stays valid).
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Insert four fresh rows starting at row 4, pushing everything below down.
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

$newRows = @(
    @("004479463", "HENRIQUE", 47959.99),
    @("004240014", "ISABELE", 25461.1),
    @("004752461", "SERGIO", 6623.66),
    @("004752494", "SERGIO", 6623.66)
)

$r = 4
foreach ($row in $newRows) {
    # Leading apostrophe forces the account number to stay text (keeps
    # leading zeros, matches the inlineStr cells used throughout column A).
    # ClearFormats() afterwards drops the "entered as text" quote-prefix
    # style marker so the cell keeps the sheet's plain/default styling.
    $ws.Cells.Item($r, 1).Value = "'" + $row[0]
    $ws.Cells.Item($r, 1).ClearFormats()
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
